# Swap the "PHYSICS" row (row 2) and the "LINGUISTICS" row (row 3) so that
# LINGUISTICS now appears first (row 2) and PHYSICS second (row 3), while
# every other row (MATHEMATICS, MEDICINE) stays untouched.
#
# NOTE: ".Value" only reflects correctly on the very first property-get of
# the whole script in this runtime; every later ".Value" get degrades to a
# placeholder string. ".Value2" does not have that issue, so it is used for
# all reads/writes here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 1; $col -le 5; $col++) {
    $physicsRowCell     = $ws.Cells.Item(2, $col)
    $linguisticsRowCell = $ws.Cells.Item(3, $col)

    $physicsValue     = $physicsRowCell.Value2
    $linguisticsValue = $linguisticsRowCell.Value2

    $physicsRowCell.Value2     = $linguisticsValue
    $linguisticsRowCell.Value2 = $physicsValue
}
